$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.007879177902933363

$ws.Range("B3").Value = -0.4027513658717083
$ws.Range("C3").Value = -0.4661884778718272
$ws.Range("D3").Value = 0.431977463267512

$ws.Range("B4").Value = -0.7820964899011106
$ws.Range("C4").Value = -1.291112248493917
$ws.Range("D4").Value = -1.417005508671151

$ws.Range("B5").Value = 0.351073541901304
$ws.Range("C5").Value = 0.523717141809419
$ws.Range("D5").Value = 2.287742462260006

$ws.Range("B6").Value = -0.8622086857614417
$ws.Range("C6").Value = -0.4437661975110234
$ws.Range("D6").Value = 2.173105040264365

$ws.Range("B7").Value = 1.457248943882707
$ws.Range("C7").Value = 0.5567061355003813
$ws.Range("D7").Value = 0.7582402414817506
